$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the first
#    (Heading1) paragraph: an empty run, a bold "Meta description" run, and
#    a plain run with the rest of the sentence.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Join the Maori gods on an adventure in Brood of Gods. Play this unique slot game for free and enjoy stunning graphics and exciting features.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2. Near the end of the document, remove the duplicated bold
#    "Play Brood of Gods Free: Unique Layout, Animations!" paragraph (the
#    real Heading1 title at the top of the doc must stay untouched, so the
#    search starts after the meta-description paragraph we just inserted).
# ---------------------------------------------------------------------------
$searchStart = $d.Paragraphs.Item(3).Range.Start
$searchRange = $d.Range($searchStart, $d.Content.End)
$found = $searchRange.Find.Execute("Play Brood of Gods Free: Unique Layout, Animations!",
                                    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Find collapses searchRange to the matched text; extend by one char to
    # also sweep up the paragraph mark so the whole paragraph disappears.
    $boldParaRange = $d.Range($searchRange.Start, $searchRange.End + 1)
    $boldParaRange.Delete()
}

# ---------------------------------------------------------------------------
# 3. Replace the (now final) italic paragraph's text with the new
#    image-generation prompt copy.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$italicPara = $d.Paragraphs.Item($n)
[void]$italicPara.Range.Find.Execute(
    "Join the Maori gods on an adventure in Brood of Gods. Play this unique slot game for free and enjoy stunning graphics and exciting features.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Design a Cartoon Style Feature Image for Brood of Gods with a Happy Maya Warrior with Glasses Create an eye-catching feature image for Brood of Gods that captures the spirit of the game. The image should be in a cartoon style, featuring a happy Maya warrior with glasses. The warrior should be surrounded by the Maori gods and goddesses featured in the game, including Papa, Rangi, and Whiro. The image should be colorful and vibrant, with a background that evokes the lush green forest where the game is set. The overall feel should be fun and adventurous, inviting players to explore the world of Brood of Gods.",
    2)

Write-Host "Done"
